$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 418-419, pushing existing rows 418..433 down to 420..435.
$ws.Range("A418:A419").EntireRow.Insert()

# Row 418 - new weekly entry (Primera)
$ws.Cells.Item(418,1).Value2  = 1
$ws.Cells.Item(418,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(418,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(418,4).Value2  = 45008
$ws.Cells.Item(418,5).Value2  = 15
$ws.Cells.Item(418,6).Value2  = 100114014
$ws.Cells.Item(418,7).Value2  = "Betarraga"
$ws.Cells.Item(418,8).Value2  = "Sin especificar"
$ws.Cells.Item(418,9).Value2  = "Primera"
$ws.Cells.Item(418,10).Value2 = 500
$ws.Cells.Item(418,11).Value2 = 700
$ws.Cells.Item(418,12).Value2 = 800
$ws.Cells.Item(418,13).Value2 = 740
$ws.Cells.Item(418,14).Value2 = "$/paquete 4 unidades"
$ws.Cells.Item(418,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(418,16).Value2 = 185
$ws.Cells.Item(418,17).Value2 = 4
$ws.Cells.Item(418,18).Value2 = "Hortaliza"

# Row 419 - new weekly entry (Segunda)
$ws.Cells.Item(419,1).Value2  = 1
$ws.Cells.Item(419,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(419,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(419,4).Value2  = 45008
$ws.Cells.Item(419,5).Value2  = 15
$ws.Cells.Item(419,6).Value2  = 100114014
$ws.Cells.Item(419,7).Value2  = "Betarraga"
$ws.Cells.Item(419,8).Value2  = "Sin especificar"
$ws.Cells.Item(419,9).Value2  = "Segunda"
$ws.Cells.Item(419,10).Value2 = 750
$ws.Cells.Item(419,11).Value2 = 700
$ws.Cells.Item(419,12).Value2 = 800
$ws.Cells.Item(419,13).Value2 = 747
$ws.Cells.Item(419,14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(419,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(419,16).Value2 = 149
$ws.Cells.Item(419,17).Value2 = 5
$ws.Cells.Item(419,18).Value2 = "Hortaliza"
